# Alloy_Steel_Flanged_Button_Head_Screws.xlsx edit
# Insert two new rows at the top of the sheet:
#   - new row 1: numeric column-index header (0..12), using the original header style
#   - new row 2: a mostly blank row with "Drive" in column E
# The previous header row (row 1) and all data rows shift down by two rows,
# and the (now relocated) header row loses its bold/border styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M")

# Insert two blank rows above the current row 1 (shifts everything, including
# the old header row and all data rows, down by two).
$ws.Range("A1:A2").EntireRow.Insert()

# The former header row (bold font + border + center/top alignment) is now row 3.
# Copy that formatting up to the new row 1, which will hold the numeric headers.
$ws.Range("A3:M3").Copy()
$ws.Range("A1:M1").PasteSpecial(-4122)   # xlPasteFormats

# Clear the (now misplaced) header formatting from row 3 - it becomes a plain
# text row again, matching the rest of the data rows.
$ws.Range("A3:M3").ClearFormats()

# Write the new row 1 values: simple 0-based column indices.
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $i
}

# Touch every cell in row 2 so each column has a (blank) cell present, then
# clear the formatting used to "touch" them so they stay unstyled.
foreach ($col in $cols) {
    $ws.Range($col + "2").NumberFormat = "@"
}
$ws.Range("E2").Value = "Drive"
$ws.Range("A2:M2").ClearFormats()

# Keep the initial selection on A1.
$ws.Range("A1").Select() | Out-Null
